$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Overview": add two new handoff-ready rows (3bf7a7e2, 99666f0c)
# ahead of the existing b69fb0f1 row, which shifts from row 3 to row 5.
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(3).Copy()
$wsOverview.Rows.Item(3).Insert()
$wsOverview.Rows.Item(3).Copy()
$wsOverview.Rows.Item(3).Insert()

$wsOverview.Range("A3").Value = "3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md"
$wsOverview.Range("B3").Value = "e2e\3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 18:45:26"

$wsOverview.Range("A4").Value = "99666f0c-86cd-436d-8323-4792d84bad10.md"
$wsOverview.Range("B4").Value = "e2e\99666f0c-86cd-436d-8323-4792d84bad10.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-02 18:45:26"

$wsOverview.Range("A5").Value = "b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md"
$wsOverview.Range("B5").Value = "e2e\b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-02 18:43:58"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3c9f54b62b3affbbaea2af1c885d7e948d50228/e2e/3aa705ce-3840-4a28-8a72-95576b742e45.md", "", "", "e2e\3aa705ce-3840-4a28-8a72-95576b742e45.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de88ab303d84295f2f04c974b2802f2c65c9d9f0/e2e/3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md", "", "", "e2e\3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c176ec483d146fbd0d741f247b65a068d124876d/e2e/99666f0c-86cd-436d-8323-4792d84bad10.md", "", "", "e2e\99666f0c-86cd-436d-8323-4792d84bad10.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/35580716ed04e5a733ddaf2555a7281d86ab90bb/e2e/b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md", "", "", "e2e\b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md") | Out-Null

# ------------------------------------------------------------------
# Sheet "zh-cn": same two new rows, localized-status columns (A..P)
# ------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Rows.Item(3).Copy()
$wsZhCn.Rows.Item(3).Insert()
$wsZhCn.Rows.Item(3).Copy()
$wsZhCn.Rows.Item(3).Insert()

$wsZhCn.Range("A3").Value = "3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "3bf7a7e2-4cec-4a8a-b147-faeca04b5900.d8fb3559cec26fdeb3853a97a0a9845a65edb402.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-02 18:45:21"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Range("A4").Value = "99666f0c-86cd-436d-8323-4792d84bad10.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "99666f0c-86cd-436d-8323-4792d84bad10.817b49e97d6fe963e10e19e5cf14069e2e8e5c75.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-02 18:45:21"
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$wsZhCn.Range("A5").Value = "b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "b69fb0f1-7d00-407b-828d-1cf89a38e6f4.1b4b06bbfb47ceeeb439a8eaec606aadf11b8197.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-09-02 18:43:54"
$wsZhCn.Range("I5").Value = ""
$wsZhCn.Range("J5").Value = ""
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L5").Value = ""
$wsZhCn.Range("M5").Value = "True"
$wsZhCn.Range("N5").Value = ""
$wsZhCn.Range("O5").Value = "False"
$wsZhCn.Range("P5").Value = ""

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P5"))

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3c9f54b62b3affbbaea2af1c885d7e948d50228/e2e/3aa705ce-3840-4a28-8a72-95576b742e45.md", "", "", "3aa705ce-3840-4a28-8a72-95576b742e45.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/23419d90b988933b8aad9e1d9214043a80f08605/e2e/3aa705ce-3840-4a28-8a72-95576b742e45.md", "", "", "3aa705ce-3840-4a28-8a72-95576b742e45.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de88ab303d84295f2f04c974b2802f2c65c9d9f0/e2e/3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md", "", "", "3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c176ec483d146fbd0d741f247b65a068d124876d/e2e/99666f0c-86cd-436d-8323-4792d84bad10.md", "", "", "99666f0c-86cd-436d-8323-4792d84bad10.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/35580716ed04e5a733ddaf2555a7281d86ab90bb/e2e/b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md", "", "", "b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md") | Out-Null

# ------------------------------------------------------------------
# Sheet "de-de": same two new rows, localized-status columns (A..P)
# ------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Rows.Item(3).Copy()
$wsDeDe.Rows.Item(3).Insert()
$wsDeDe.Rows.Item(3).Copy()
$wsDeDe.Rows.Item(3).Insert()

$wsDeDe.Range("A3").Value = "3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "3bf7a7e2-4cec-4a8a-b147-faeca04b5900.d8fb3559cec26fdeb3853a97a0a9845a65edb402.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-02 18:45:26"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Range("A4").Value = "99666f0c-86cd-436d-8323-4792d84bad10.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "99666f0c-86cd-436d-8323-4792d84bad10.817b49e97d6fe963e10e19e5cf14069e2e8e5c75.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-02 18:45:26"
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$wsDeDe.Range("A5").Value = "b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "b69fb0f1-7d00-407b-828d-1cf89a38e6f4.1b4b06bbfb47ceeeb439a8eaec606aadf11b8197.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-09-02 18:43:58"
$wsDeDe.Range("I5").Value = ""
$wsDeDe.Range("J5").Value = ""
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L5").Value = ""
$wsDeDe.Range("M5").Value = "True"
$wsDeDe.Range("N5").Value = ""
$wsDeDe.Range("O5").Value = "False"
$wsDeDe.Range("P5").Value = ""

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P5"))

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3c9f54b62b3affbbaea2af1c885d7e948d50228/e2e/3aa705ce-3840-4a28-8a72-95576b742e45.md", "", "", "3aa705ce-3840-4a28-8a72-95576b742e45.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8734e807321c985e461554f0f835d2e1240cec8c/e2e/3aa705ce-3840-4a28-8a72-95576b742e45.md", "", "", "3aa705ce-3840-4a28-8a72-95576b742e45.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de88ab303d84295f2f04c974b2802f2c65c9d9f0/e2e/3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md", "", "", "3bf7a7e2-4cec-4a8a-b147-faeca04b5900.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c176ec483d146fbd0d741f247b65a068d124876d/e2e/99666f0c-86cd-436d-8323-4792d84bad10.md", "", "", "99666f0c-86cd-436d-8323-4792d84bad10.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/35580716ed04e5a733ddaf2555a7281d86ab90bb/e2e/b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md", "", "", "b69fb0f1-7d00-407b-828d-1cf89a38e6f4.md") | Out-Null

"Generate Report for Handoff: done"
